# Refreshes the cryptos price list (Price / Volume(1h) columns, plus the
# Elrond <-> EnergySwap row swap at rows 47/48) to match the scraped update.
#
# These cells hold literal text (e.g. "1.000", "30.536.42", "  +0.64%  ")
# rather than numbers, so every new value below is written with a leading
# apostrophe -- Excel's standard "force text" marker -- to stop values like
# "1.000" or "18.00" from being auto-coerced into numbers on assignment.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$apos = "'"   # leading apostrophe = force-text prefix (kept as a real string)

$ws.Range("D2").Value = $apos + '30.536.42'
$ws.Range("E2").Value = $apos + '  +0.64%  '
$ws.Range("D3").Value = $apos + '1.852.76'
$ws.Range("E3").Value = $apos + '  -0.35%  '
$ws.Range("D4").Value = $apos + '1.000'
$ws.Range("E4").Value = $apos + '  -0.05%  '
$ws.Range("D5").Value = $apos + '233.60'
$ws.Range("E5").Value = $apos + '  -0.41%  '
$ws.Range("D7").Value = $apos + '0.4718'
$ws.Range("E7").Value = $apos + '  -0.42%  '
$ws.Range("D8").Value = $apos + '0.2743'
$ws.Range("E8").Value = $apos + '  -0.01%  '
$ws.Range("D9").Value = $apos + '0.06328'
$ws.Range("E9").Value = $apos + '  -1.59%  '
$ws.Range("D10").Value = $apos + '17.66'
$ws.Range("E10").Value = $apos + '  +8.12%  '
$ws.Range("D11").Value = $apos + '1.820.17'
$ws.Range("E11").Value = $apos + '  -0.62%  '
$ws.Range("D12").Value = $apos + '0.07406'
$ws.Range("E12").Value = $apos + '  -0.71%  '
$ws.Range("D13").Value = $apos + '5.038'
$ws.Range("E13").Value = $apos + '  +0.92%  '
$ws.Range("D14").Value = $apos + '84.51'
$ws.Range("E14").Value = $apos + '  -1.18%  '
$ws.Range("D15").Value = $apos + '0.6251'
$ws.Range("E15").Value = $apos + '  -1.25%  '
$ws.Range("D16").Value = $apos + '30.490.22'
$ws.Range("E16").Value = $apos + '  +0.61%  '
$ws.Range("D17").Value = $apos + '242.64'
$ws.Range("E17").Value = $apos + '  +5.46%  '
$ws.Range("D18").Value = $apos + '1.001'
$ws.Range("E18").Value = $apos + '  -0.02%  '
$ws.Range("D19").Value = $apos + '12.64'
$ws.Range("E19").Value = $apos + '  -0.81%  '
$ws.Range("D20").Value = $apos + '0.000007345'
$ws.Range("E20").Value = $apos + '  -0.74%  '
$ws.Range("E21").Value = $apos + '  +0.11%  '
$ws.Range("D22").Value = $apos + '4.930'
$ws.Range("E22").Value = $apos + '  -1.19%  '
$ws.Range("D23").Value = $apos + '5.963'
$ws.Range("E23").Value = $apos + '  -0.53%  '
$ws.Range("D24").Value = $apos + '9.212'
$ws.Range("E24").Value = $apos + '  -0.76%  '
$ws.Range("D25").Value = $apos + '161.31'
$ws.Range("E25").Value = $apos + '  -3.13%  '
$ws.Range("D26").Value = $apos + '18.00'
$ws.Range("E26").Value = $apos + '  +0.14%  '
$ws.Range("D27").Value = $apos + '1.881'
$ws.Range("E27").Value = $apos + '  -0.46%  '
$ws.Range("D28").Value = $apos + '0.1020'
$ws.Range("E28").Value = $apos + '  -2.92%  '
$ws.Range("E29").Value = $apos + '  -2.70%  '
$ws.Range("D30").Value = $apos + '4.019'
$ws.Range("E30").Value = $apos + '  -3.16%  '
$ws.Range("D31").Value = $apos + '3.832'
$ws.Range("E31").Value = $apos + '  -2.46%  '
$ws.Range("D32").Value = $apos + '0.04861'
$ws.Range("E32").Value = $apos + '  -1.56%  '
$ws.Range("D33").Value = $apos + '1.135'
$ws.Range("E33").Value = $apos + '  -2.63%  '
$ws.Range("D34").Value = $apos + '0.7056'
$ws.Range("E34").Value = $apos + '  -2.67%  '
$ws.Range("D35").Value = $apos + '2.714'
$ws.Range("E35").Value = $apos + '  +0.50%  '
$ws.Range("D36").Value = $apos + '0.01906'
$ws.Range("E36").Value = $apos + '  +1.83%  '
$ws.Range("D37").Value = $apos + '2.683'
$ws.Range("E37").Value = $apos + '  +1.40%  '
$ws.Range("D38").Value = $apos + '0.8746'
$ws.Range("E38").Value = $apos + '  -4.38%  '
$ws.Range("D39").Value = $apos + '1.974'
$ws.Range("E39").Value = $apos + '  +0.09%  '
$ws.Range("D40").Value = $apos + '105.32'
$ws.Range("E40").Value = $apos + '  -0.79%  '
$ws.Range("D41").Value = $apos + '1.001'
$ws.Range("E41").Value = $apos + '  +0.06%  '
$ws.Range("D42").Value = $apos + '0.4070'
$ws.Range("E42").Value = $apos + '  -1.03%  '
$ws.Range("D43").Value = $apos + '5.500'
$ws.Range("E43").Value = $apos + '  -1.12%  '
$ws.Range("D44").Value = $apos + '7.208'
$ws.Range("E44").Value = $apos + '  +1.15%  '
$ws.Range("D45").Value = $apos + '62.26'
$ws.Range("E45").Value = $apos + '  +2.25%  '
$ws.Range("D46").Value = $apos + '0.1213'
$ws.Range("E46").Value = $apos + '  +1.26%  '
$ws.Range("B47").Value = $apos + 'EnergySwap'
$ws.Range("C47").Value = $apos + 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").Value = $apos + '8.539'
$ws.Range("E47").Value = $apos + '  -1.84%  '
$ws.Range("B48").Value = $apos + 'Elrond'
$ws.Range("C48").Value = $apos + 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D48").Value = $apos + '33.30'
$ws.Range("E48").Value = $apos + '  -0.44%  '
$ws.Range("D49").Value = $apos + '0.05543'
$ws.Range("E49").Value = $apos + '  -0.72%  '
$ws.Range("D50").Value = $apos + '1.367'
$ws.Range("E50").Value = $apos + '  -2.88%  '
$ws.Range("D51").Value = $apos + '0.3669'
$ws.Range("E51").Value = $apos + '  -0.87%  '
